$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.013.40'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '3.768.56'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '403.78'
$ws.Range('E5').Value = '  -4.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.35'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').Value = '3.759.97'
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').Value = '  -6.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -6.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.169'
$ws.Range('E11').Value = '  -9.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000360'
$ws.Range('E12').Value = '  -12.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.65'
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').Value = '4.365.69'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.71'
$ws.Range('E15').Value = '  -4.85%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.63'
$ws.Range('E16').Value = '  +12.54%  '
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '3.771.44'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.43'
$ws.Range('E19').Value = '  -6.84%  '
$ws.Range('D20').Value = '66.204.89'
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('E21').Value = '  -6.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '410.95'
$ws.Range('E22').Value = '  -9.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.39'
$ws.Range('E23').Value = '  -7.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.82'
$ws.Range('E24').Value = '  -5.61%  '
$ws.Range('E25').Value = '  -4.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.74'
$ws.Range('E26').Value = '  +15.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '36.05'
$ws.Range('E27').Value = '  -5.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.10'
$ws.Range('E28').Value = '  -7.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.34'
$ws.Range('E29').Value = '  -8.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.31'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.118'
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.38'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  -6.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '38.97'
$ws.Range('E35').Value = '  -7.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.02'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').Value = '0.0₃0735'
$ws.Range('E38').Value = '  -7.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0458'
$ws.Range('E39').Value = '  -7.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.89'
$ws.Range('E40').Value = '  -6.66%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  -8.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '27.16'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '144.77'
$ws.Range('E44').Value = '  -1.84%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.09'
$ws.Range('E45').Value = '  +16.90%  '
$ws.Range('E46').Value = '  -4.91%  '
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.57'
$ws.Range('E48').Value = '  -4.28%  '
$ws.Range('E49').Value = '  -5.66%  '
$ws.Range('E50').Value = '  -5.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.293'
$ws.Range('E51').Value = '  -4.89%  '
